$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style/formatting (bold, border, centered) from the existing
# header cell AC1 onto the three new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every player row (2-52)
$ws.Range("AD2:AD52").Value = 86
$ws.Range("AE2:AE52").Value = 76
$ws.Range("AF2:AF52").Value = 0
